$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.748.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.437.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.33%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.31%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("E9").Value = "  +6.67%  "

$ws.Range("E10").Value = "  -2.09%  "

$ws.Range("E11").Value = "  -0.80%  "

$ws.Range("E12").Value = "  -5.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000176"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "68.655.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.885.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.437.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.93%  "

$ws.Range("E22").Value = "  +1.87%  "

$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("E25").Value = "  +1.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.563.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.21"
$ws.Range("D28").Style = "Normal"

$ws.Range("E29").Value = "  -0.82%  "

$ws.Range("E30").Value = "  -1.28%  "

$ws.Range("E31").Value = "  +0.03%  "

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Bittensor"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "427.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.60%  "

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.02%  "

$ws.Range("E34").Value = "  -2.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.17%  "

$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("E39").Value = "  -2.83%  "

$ws.Range("E40").Value = "  -0.29%  "

$ws.Range("E41").Value = "  +2.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.38%  "

$ws.Range("E43").Value = "  -0.75%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Aave"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "130.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.42%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Filecoin"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.36%  "

$ws.Range("E47").Value = "  -0.19%  "

$ws.Range("E48").Value = "  -0.87%  "

$ws.Range("E49").Value = "  -1.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0921"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.42%  "

$ws.Range("E51").Value = "  +2.19%  "
